# Update UV_vis_Template header row (row 1) to the new "inv template" schema
# and shrink the used range from A1:AE1 down to A1:Q1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Absorber Emitter_inchikey <absorber_emitter>",
    "Absorber Emitter_molfile",
    "Experiment Type <experiment_type>",
    "Absorption Max [nm] <absorption_max>",
    "Absorption Intensity Max [nm] <absorption_intensity_max>",
    "Emission Max [nm] <emission_max>",
    "Emission Intensity Max [nm] <emission_intensity_max>",
    "Intersection [nm] <intersection>",
    "Solvent_inchikey <solvent>",
    "Solvent_molfile",
    "Solvent concentration [µM] <solvent_conc>",
    "TemperatureP [°C] <temp>",
    "E0 [eV] <auto-generated-E0>",
    "Absorption Coefficient [L/mol⋅cm] <absorption_coefficient>",
    "Details <details>",
    "Included <include>",
    "BasePageName <BasePageName>"
)

# Overwrite A1:Q1 in place with the new labels (keeps existing yellow header style).
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Drop the now-unused trailing columns (old sheet went out to AE1 / 31 columns,
# new one only needs 17 -> A:Q), shifting the used range/dimension down to A1:Q1.
$ws.Range("R1:AE1").EntireColumn.Delete()

# Widen column M (13) to fit the longer "E0 [eV] <auto-generated-E0>" label.
$ws.Columns.Item(13).ColumnWidth = 25.85
